$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3798.9312
$ws.Cells.Item(40, 9).Value = 2495
$ws.Cells.Item(40, 11).Value = 2495
$ws.Cells.Item(40, 13).Value = -2320
$ws.Cells.Item(69, 8).Value = 6998
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 14).ClearContents()
$ws.Cells.Item(72, 8).Value = 6998
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 14).ClearContents()
$ws.Cells.Item(137, 8).Value = 41845.668
$ws.Cells.Item(137, 9).Value = 84567.39999999999
$ws.Cells.Item(137, 11).Value = 253702.2
$ws.Cells.Item(137, 13).Value = -251152.2
$ws.Cells.Item(138, 8).Value = 2471.5715
$ws.Cells.Item(138, 9).Value = 1565.4
$ws.Cells.Item(138, 11).Value = 4696.200000000001
$ws.Cells.Item(138, 13).Value = 443.7999999999993
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 3117.6538
$ws.Cells.Item(74, 9).Value = 2823.2917
$ws.Cells.Item(74, 11).Value = 2823.2917
$ws.Cells.Item(74, 13).Value = -1949.2917
$ws.Cells.Item(77, 8).Value = 3117.6538
$ws.Cells.Item(77, 9).Value = 2823.2917
$ws.Cells.Item(77, 11).Value = 14116.4585
$ws.Cells.Item(77, 13).Value = -9748.458500000001
$ws.Cells.Item(92, 8).Value = 67994.5
$ws.Cells.Item(92, 10).Value = 67994.5
$ws.Cells.Item(92, 12).Value = 67994.5
$ws.Cells.Item(92, 14).Value = -72986.5
$ws.Cells.Item(102, 8).Value = 1950.3125
$ws.Cells.Item(102, 10).Value = 3749.5
$ws.Cells.Item(102, 12).Value = 3749.5
$ws.Cells.Item(102, 14).Value = -6993.5
$ws.Cells.Item(110, 8).Value = 1908.3
$ws.Cells.Item(110, 9).Value = 1761
$ws.Cells.Item(110, 10).Value = 2497.5
$ws.Cells.Item(110, 11).Value = 1761
$ws.Cells.Item(110, 12).Value = 2497.5
$ws.Cells.Item(110, 13).Value = 284
$ws.Cells.Item(110, 14).Value = -6587.5
$ws.Cells.Item(122, 8).Value = 2978.4243
$ws.Cells.Item(122, 9).Value = 2407.7827
$ws.Cells.Item(122, 11).Value = 7223.348100000001
$ws.Cells.Item(122, 13).Value = -4773.348100000001
$ws.Cells.Item(132, 8).Value = 265938.9
$ws.Cells.Item(132, 9).Value = 347109.22
$ws.Cells.Item(132, 10).Value = 4390.222
$ws.Cells.Item(132, 11).Value = 1041327.66
$ws.Cells.Item(132, 12).Value = 13170.666
$ws.Cells.Item(132, 13).Value = -1038797.66
$ws.Cells.Item(132, 14).Value = -18230.666
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2854.8096
$ws.Cells.Item(105, 9).Value = 2546.4167
$ws.Cells.Item(105, 11).Value = 2546.4167
$ws.Cells.Item(105, 13).Value = -799.4167000000002
$ws.Cells.Item(116, 8).Value = 83408
$ws.Cells.Item(116, 10).Value = 83408
$ws.Cells.Item(116, 12).Value = 83408
$ws.Cells.Item(116, 14).Value = -92586
$ws.Cells.Item(134, 8).Value = 3405333
$ws.Cells.Item(134, 9).Value = 5104221
$ws.Cells.Item(134, 11).Value = 15312663
$ws.Cells.Item(134, 13).Value = -15310128
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(18, 8).Value = 25864.5
$ws.Cells.Item(18, 10).Value = 25864.5
$ws.Cells.Item(18, 12).Value = 25864.5
$ws.Cells.Item(18, 14).Value = -26324.5
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 13).ClearContents()
$ws.Cells.Item(31, 8).Value = 7188.3335
$ws.Cells.Item(31, 9).Value = 3093.0715
$ws.Cells.Item(31, 11).Value = 3093.0715
$ws.Cells.Item(31, 13).Value = -2798.0715
$ws.Cells.Item(34, 8).Value = 7188.3335
$ws.Cells.Item(34, 9).Value = 3093.0715
$ws.Cells.Item(34, 11).Value = 3093.0715
$ws.Cells.Item(34, 13).Value = -2891.0715
$ws.Cells.Item(36, 8).Value = 40009.25
$ws.Cells.Item(36, 9).Value = 40009.25
$ws.Cells.Item(36, 11).Value = 40009.25
$ws.Cells.Item(36, 13).Value = -39621.25
$ws.Cells.Item(40, 8).Value = 40009.25
$ws.Cells.Item(40, 9).Value = 40009.25
$ws.Cells.Item(40, 11).Value = 40009.25
$ws.Cells.Item(40, 13).Value = -39849.25
$ws.Cells.Item(42, 8).Value = 1000
$ws.Cells.Item(42, 9).Value = 1000
$ws.Cells.Item(42, 11).Value = 1000
$ws.Cells.Item(42, 13).Value = -407
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 13).ClearContents()
$ws.Cells.Item(52, 8).Value = 32799
$ws.Cells.Item(52, 10).Value = 32799
$ws.Cells.Item(52, 12).Value = 32799
$ws.Cells.Item(52, 14).Value = -33387
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(62, 8).Value = 4950
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 4950
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 4950
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 14).Value = -6198
$ws.Cells.Item(65, 8).Value = 4950
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 4950
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 24750
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 14).Value = -30990
$ws.Cells.Item(69, 8).Value = 54659.75
$ws.Cells.Item(69, 9).Value = 40399.8
$ws.Cells.Item(69, 10).Value = 78426.336
$ws.Cells.Item(69, 11).Value = 40399.8
$ws.Cells.Item(69, 12).Value = 78426.336
$ws.Cells.Item(69, 13).Value = -39650.8
$ws.Cells.Item(69, 14).Value = -79924.336
$ws.Cells.Item(72, 8).Value = 54659.75
$ws.Cells.Item(72, 9).Value = 40399.8
$ws.Cells.Item(72, 10).Value = 78426.336
$ws.Cells.Item(72, 11).Value = 121199.4
$ws.Cells.Item(72, 12).Value = 235279.008
$ws.Cells.Item(72, 13).Value = -117455.4
$ws.Cells.Item(72, 14).Value = -242767.008
$ws.Cells.Item(99, 8).Value = 2858.7856
$ws.Cells.Item(99, 10).Value = 3199.4
$ws.Cells.Item(99, 12).Value = 3199.4
$ws.Cells.Item(99, 14).Value = -6195.4
$ws.Cells.Item(117, 8).Value = 57249.25
$ws.Cells.Item(117, 10).Value = 57249.25
$ws.Cells.Item(117, 12).Value = 57249.25
$ws.Cells.Item(117, 14).Value = -66427.25
$ws.Cells.Item(118, 8).Value = 115995
$ws.Cells.Item(118, 10).Value = 112990
$ws.Cells.Item(118, 12).Value = 112990
$ws.Cells.Item(118, 14).Value = -116304
$ws.Cells.Item(126, 8).Value = 2858.7856
$ws.Cells.Item(126, 10).Value = 3199.4
$ws.Cells.Item(126, 12).Value = 9598.200000000001
$ws.Cells.Item(126, 14).Value = -14538.2
$ws.Cells.Item(132, 8).Value = 3555.9375
$ws.Cells.Item(132, 9).Value = 3555.9375
$ws.Cells.Item(132, 11).Value = 10667.8125
$ws.Cells.Item(132, 13).Value = -8137.8125
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 2.6
$ws.Cells.Item(12, 10).Value = 2.6
$ws.Cells.Item(12, 12).Value = 7.800000000000001
$ws.Cells.Item(12, 14).Value = -353.8
$ws.Cells.Item(117, 8).Value = 4655.1
$ws.Cells.Item(117, 10).Value = 4812
$ws.Cells.Item(117, 12).Value = 14436
$ws.Cells.Item(117, 14).Value = -21320
$ws.Cells.Item(132, 8).Value = 932.2778
$ws.Cells.Item(132, 9).Value = 975
$ws.Cells.Item(132, 10).Value = 926.9375
$ws.Cells.Item(132, 11).Value = 8775
$ws.Cells.Item(132, 12).Value = 8342.4375
$ws.Cells.Item(132, 13).Value = -6245
$ws.Cells.Item(132, 14).Value = -13402.4375
$ws.Cells.Item(134, 8).Value = 6719.1875
$ws.Cells.Item(134, 9).Value = 6137.273
$ws.Cells.Item(134, 11).Value = 18411.819
$ws.Cells.Item(134, 13).Value = -13341.819
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 4414499.5
$ws.Cells.Item(11, 10).Value = 1583285.8
$ws.Cells.Item(11, 12).Value = 1583285.8
$ws.Cells.Item(11, 14).Value = -1583563.8
$ws.Cells.Item(97, 8).Value = 842.7059
$ws.Cells.Item(97, 9).Value = 537.7273
$ws.Cells.Item(97, 11).Value = 537.7273
$ws.Cells.Item(97, 13).Value = -41.72730000000001
$ws.Cells.Item(116, 8).Value = 106889
$ws.Cells.Item(116, 10).Value = 106889
$ws.Cells.Item(116, 12).Value = 106889
$ws.Cells.Item(116, 14).Value = -116067
$ws.Cells.Item(122, 8).Value = 1241.75
$ws.Cells.Item(122, 9).Value = 1241.75
$ws.Cells.Item(122, 11).Value = 3725.25
$ws.Cells.Item(122, 13).Value = -1275.25
$ws.Cells.Item(124, 8).Value = 152599
$ws.Cells.Item(124, 10).Value = 152599
$ws.Cells.Item(124, 12).Value = 152599
$ws.Cells.Item(124, 14).Value = -162419
$ws.Cells.Item(132, 8).Value = 5166.5557
$ws.Cells.Item(132, 9).Value = 5253.6924
$ws.Cells.Item(132, 10).Value = 4940
$ws.Cells.Item(132, 11).Value = 15761.0772
$ws.Cells.Item(132, 12).Value = 14820
$ws.Cells.Item(132, 13).Value = -13231.0772
$ws.Cells.Item(132, 14).Value = -19880
$ws.Cells.Item(137, 8).Value = 46333.332
$ws.Cells.Item(137, 10).Value = 46785.715
$ws.Cells.Item(137, 12).Value = 46785.715
$ws.Cells.Item(137, 14).Value = -56985.715
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4600.7
$ws.Cells.Item(40, 9).Value = 2000.2858
$ws.Cells.Item(40, 10).Value = 10668.333
$ws.Cells.Item(40, 11).Value = 2000.2858
$ws.Cells.Item(40, 12).Value = 10668.333
$ws.Cells.Item(40, 13).Value = -1864.2858
$ws.Cells.Item(40, 14).Value = -10940.333
$ws.Cells.Item(81, 8).Value = 114989
$ws.Cells.Item(81, 10).Value = 114989
$ws.Cells.Item(81, 12).Value = 114989
$ws.Cells.Item(81, 14).Value = -116985
$ws.Cells.Item(84, 8).Value = 114989
$ws.Cells.Item(84, 10).Value = 114989
$ws.Cells.Item(84, 12).Value = 344967
$ws.Cells.Item(84, 14).Value = -354951
$ws.Cells.Item(132, 8).Value = 358128.9
$ws.Cells.Item(132, 9).Value = 442808.53
$ws.Cells.Item(132, 11).Value = 1328425.59
$ws.Cells.Item(132, 13).Value = -1325895.59
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 102478.75
$ws.Cells.Item(16, 10).Value = 102478.75
$ws.Cells.Item(16, 12).Value = 102478.75
$ws.Cells.Item(16, 14).Value = -103062.75
$ws.Cells.Item(131, 8).Value = 134000
$ws.Cells.Item(131, 10).Value = 134000
$ws.Cells.Item(131, 12).Value = 134000
$ws.Cells.Item(131, 14).Value = -144080
$ws.Cells.Item(132, 8).Value = 23751.305
$ws.Cells.Item(132, 9).Value = 25682.762
$ws.Cells.Item(132, 11).Value = 77048.28599999999
$ws.Cells.Item(132, 13).Value = -74518.28599999999
$ws.Cells.Item(136, 8).Value = 45911
$ws.Cells.Item(136, 9).Value = 1672.8
$ws.Cells.Item(136, 11).Value = 5018.4
$ws.Cells.Item(136, 13).Value = -2468.4
